$A1 = 'TabName'
$B1 = 'query'
$C1 = 'StatQuery'
$D1 = 'dbExcel'
$E1 = 'WebExcel'
$A2 = 'CasesTab'
$B2 = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
 MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
WHERE diag.best_response IN ["Complete Response"] 
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@
$C2 = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (cf:file)-->(samp)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT cf, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.best_response IN ["Complete Response"]
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$D2 = 'TC01_Canine_Filter_RespToTrtmt-CompleteResponse_Neo4jData.xlsx'
$E2 = 'TC01_Canine_Filter_RespToTrtmt-CompleteResponse_WebData.xlsx'
$A3 = 'SamplesTab'
$B3 = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
MATCH (f:file)-[*]->(c)
MATCH (f:file)-->(parent)
WHERE diag.best_response IN ["Complete Response"]
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 100
'@
$C3 = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (cf:file)-->(samp)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT cf, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.best_response IN ["Complete Response"]
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$D3 = 'TC01_Canine_Filter_RespToTrtmt-CompleteResponse_Neo4jData.xlsx'
$E3 = 'TC01_Canine_Filter_RespToTrtmt-CompleteResponse_WebData.xlsx'
$A4 = 'FilesTab'
$B4 = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
optional MATCH (f)-->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.best_response IN ["Complete Response"]
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 100
'@
$C4 = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (cf:file)-->(samp)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT cf, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.best_response IN ["Complete Response"]
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$D4 = 'TC01_Canine_Filter_RespToTrtmt-CompleteResponse_Neo4jData.xlsx'
$E4 = 'TC01_Canine_Filter_RespToTrtmt-CompleteResponse_WebData.xlsx'
$A5 = 'StudyFilesTab'
$B5 = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
WHERE diag.best_response IN ["Complete Response"]
MATCH (c)<--(demo:demographic)
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@
$C5 = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (cf:file)-->(samp)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT cf, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.best_response IN ["Complete Response"]
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$D5 = 'TC01_Canine_Filter_RespToTrtmt-CompleteResponse_Neo4jData.xlsx'
$E5 = 'TC01_Canine_Filter_RespToTrtmt-CompleteResponse_WebData.xlsx'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1) - unchanged text, but now carries the new (larger)
#    font used throughout the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = $A1
$ws.Range("B1").Value = $B1
$ws.Range("C1").Value = $C1
$ws.Range("D1").Value = $D1
$ws.Range("E1").Value = $E1

# ---------------------------------------------------------------------------
# 2. Data rows 2-4 (CasesTab / SamplesTab / FilesTab) - queries rewritten,
#    stat query replaced, file names kept.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = $A2
$ws.Range("B2").Value = $B2
$ws.Range("C2").Value = $C2
$ws.Range("D2").Value = $D2
$ws.Range("E2").Value = $E2

$ws.Range("A3").Value = $A3
$ws.Range("B3").Value = $B3
$ws.Range("C3").Value = $C3
$ws.Range("D3").Value = $D3
$ws.Range("E3").Value = $E3

$ws.Range("A4").Value = $A4
$ws.Range("B4").Value = $B4
$ws.Range("C4").Value = $C4
$ws.Range("D4").Value = $D4
$ws.Range("E4").Value = $E4

# ---------------------------------------------------------------------------
# 3. New row 5 (StudyFilesTab) - brand new tab added to the table.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = $A5
$ws.Range("B5").Value = $B5
$ws.Range("C5").Value = $C5
$ws.Range("D5").Value = $D5
$ws.Range("E5").Value = $E5

# ---------------------------------------------------------------------------
# 4. Formatting - the whole table now uses a larger (16pt) font, and the
#    query / stat-query columns (B & C) are word-wrapped.
# ---------------------------------------------------------------------------
$ws.Range("A1:E5").Font.Size = 16
$ws.Range("B1:C5").WrapText = $true

# ---------------------------------------------------------------------------
# 5. Row heights - explicit custom heights as the content/font changed.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 21
$ws.Rows.Item(2).RowHeight = 237.75
$ws.Rows.Item(3).RowHeight = 188.25
$ws.Rows.Item(4).RowHeight = 150.75
$ws.Rows.Item(5).RowHeight = 409.5

# ---------------------------------------------------------------------------
# 6. Sheet view - scroll back to the top and select B4 (matches the saved
#    selection/scroll state in the edited workbook).
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("B4").Select()

Write-Host "edit complete"
